$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 30 (Item# 350 - DESICCANT PACKS BX/2 SETS), which shifts all
# subsequent rows up by one.
$ws.Rows.Item(30).Delete()

# Update the view: scroll so row 17 is the top-left visible row, and select C22.
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("C22").Select()
